$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Locate the paragraph that ends with "...varchar attributes in the
# database might be considered rather small." -- this is the last
# limitation bullet before the new ones get appended.
# ---------------------------------------------------------------------
$needle = "-The character limits for the varchar attributes in the database might be considered rather small."

$found = $d.Content
$ok = $found.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) {
    throw "Could not find the varchar-limits paragraph"
}
$matchStart = $found.Start
$matchEnd = $found.End

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $matchStart -and $p.Range.End -ge $matchEnd) {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not resolve the owning paragraph"
}

# ---------------------------------------------------------------------
# The trailing "_GoBack" bookmark currently sits at the end of that
# paragraph. Once the extra bullets below are typed, Word relocates
# _GoBack to track the newest edit point, which lands inside the new
# "Cascade deletions" paragraph (right after "...Adding the question").
# Drop it here; it is recreated further down in the new XML.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$pPr = "<w:pPr><w:pStyle w:val='segoeui'/><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr></w:pPr>"
$rPr = "<w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr>"

$noValueLimitationsText = "-No value limitations outside of those required.  Although there are some values that could have been limited, for example, keeping GPA between 0 and 4, this is not done as it was not an explicit requirement.  However, relationship constraints between the main required entities are in place, for example, a student cannot create an application for a degree that does not exist"
$cascadeDeletionsText = "-Cascade deletions.  There were not requirements on what should happen should deletion events occur, and so I kept it simple to conserve time, and have deletions cascade.  Please note that additional data may be lost when deleting something.  For example, deleting a degree requirement will delete all the answers to that question.  This is in place so a student cannot answer a question that does not exist.  Adding the question"
$cascadeDeletionsTail = " again will not recover the answers."
$trailingSpaceText = " "

$paraNoValueLimitations = "<w:p $ns>$pPr<w:r>$rPr<w:t>$noValueLimitationsText</w:t></w:r></w:p>"

$paraCascadeDeletions = "<w:p $ns>$pPr" + `
    "<w:r>$rPr<w:t>$cascadeDeletionsText</w:t></w:r>" + `
    "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" + `
    "<w:r>$rPr<w:t xml:space='preserve'>$cascadeDeletionsTail</w:t></w:r>" + `
    "</w:p>"

$emptyPara = "<w:p $ns>$pPr</w:p>"

$paraTrailingSpace = "<w:p $ns>$pPr<w:r>$rPr<w:t xml:space='preserve'>$trailingSpaceText</w:t></w:r></w:p>"

$tenEmptyParas = ""
for ($i = 0; $i -lt 10; $i++) {
    $tenEmptyParas += $emptyPara
}

$newXml = $paraNoValueLimitations + $paraCascadeDeletions + $tenEmptyParas + $paraTrailingSpace

# Insert the new bullets/paragraphs right after the varchar paragraph.
$insertAt = $d.Range($target.Range.End, $target.Range.End)
$insertAt.InsertXML($newXml)

Write-Output "OK: paragraphs now $($d.Paragraphs.Count)"
